$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.078.06'
$ws.Range('E2').Value = '  -0.76%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.469.13'
$ws.Range('E3').Value = '  -0.98%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '581.78'
$ws.Range('E5').Value = '  -1.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '167.87'
$ws.Range('E6').Value = '  -2.38%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  -1.69%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.468.27'
$ws.Range('E9').Value = '  -0.89%  '
$ws.Range('E10').Value = '  -2.12%  '
$ws.Range('E11').Value = '  -0.32%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.00'
$ws.Range('E12').Value = '  -2.11%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.331'
$ws.Range('E13').Value = '  -2.39%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.46'
$ws.Range('E14').Value = '  -3.00%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '67.021.88'
$ws.Range('E16').Value = '  -0.56%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000169'
$ws.Range('E17').Value = '  -3.75%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.525.29'
$ws.Range('E18').Value = '  +0.45%  '
$ws.Range('E19').Value = '  -4.60%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.53'
$ws.Range('E20').Value = '  -3.39%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '353.18'
$ws.Range('E21').Value = '  -3.67%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.04'
$ws.Range('E22').Value = '  -2.55%  '
$ws.Range('E23').Value = '  +0.30%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '69.13'
$ws.Range('E24').Value = '  -3.00%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.23'
$ws.Range('E25').Value = '  -7.28%  '
$ws.Range('E26').Value = '  -7.17%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.15'
$ws.Range('E27').Value = '  -7.80%  '
$ws.Range('E28').Value = '  -0.05%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.593.64'
$ws.Range('E29').Value = '  -0.32%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0903'
$ws.Range('E30').Value = '  -5.66%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '512.50'
$ws.Range('E31').Value = '  -3.50%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.71'
$ws.Range('E32').Value = '  -6.60%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.77'
$ws.Range('E33').Value = '  -4.60%  '
$ws.Range('B34').Value = 'Fetch.AI'
$ws.Range('C34').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.24'
$ws.Range('E34').Value = '  -5.63%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('E36').Value = '  -6.96%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '159.58'
$ws.Range('E37').Value = '  +1.29%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.65'
$ws.Range('E38').Value = '  +0.07%  '
$ws.Range('E39').Value = '  -2.05%  '
$ws.Range('E40').Value = '  -4.80%  '
$ws.Range('E41').Value = '  +0.01%  '
$ws.Range('B42').Value = 'PolygonEcosystemToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.327'
$ws.Range('E42').Value = '  -6.26%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.66'
$ws.Range('E43').Value = '  -5.98%  '
$ws.Range('E44').Value = '  -6.00%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.35'
$ws.Range('E45').Value = '  -4.10%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '38.66'
$ws.Range('E46').Value = '  -2.51%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '140.86'
$ws.Range('E47').Value = '  -3.02%  '
$ws.Range('E48').Value = '  -5.93%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.513'
$ws.Range('E49').Value = '  -6.16%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0₆0253'
$ws.Range('E50').Value = '  -9.22%  '
$ws.Range('B51').Value = 'Optimism'
$ws.Range('C51').Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.59'
$ws.Range('E51').Value = '  -5.69%  '
